# Insert a new daily-ranking record for 2026/02/21 07:00 (row 825), pushing
# the existing rows 825:866 down to 826:867.
#
# The sheet's "date" column stores dates as plain text (e.g. "2026/02/21"),
# not as real Excel date serials. Using Range.Value on a brand-new cell with
# a date-shaped string would be auto-converted to a date serial by Excel's
# type inference, so instead we duplicate the row immediately above (which
# already carries the correct date/weekday text with plain/default
# formatting) and then overwrite just the two numeric columns that actually
# change for the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 824 (2026/02/21, 土, 3, 78) and insert the copy above row 825,
# shifting rows 825:866 down to 826:867. The new row 825 starts out as an
# exact duplicate of row 824 (same date text "2026/02/21" and weekday "土").
$ws.Rows.Item(824).Copy()
$ws.Rows.Item(825).Insert()

# Now fix up the two numeric columns for the newly inserted row so it reads
# 2026/02/21, 土, 7, 80.
$ws.Range("C825").Value = 7
$ws.Range("D825").Value = 80
